$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 5062.625
$ws.Range("I74").Value = 4000.2307
$ws.Range("J74").Value = 6318.1816
$ws.Range("K74").Value = 4000.2307
$ws.Range("L74").Value = 6318.1816
$ws.Range("M74").Value = -3064.2307
$ws.Range("N74").Value = -8190.1816

$ws.Range("H77").Value = 5062.625
$ws.Range("I77").Value = 4000.2307
$ws.Range("J77").Value = 6318.1816
$ws.Range("K77").Value = 20001.1535
$ws.Range("L77").Value = 31590.908
$ws.Range("M77").Value = -15321.1535
$ws.Range("N77").Value = -40950.908

$ws.Range("H80").Value = 865.0714
$ws.Range("I80").Value = 814.63635
$ws.Range("J80").Value = 1050
$ws.Range("K80").Value = 2443.90905
$ws.Range("L80").Value = 3150
$ws.Range("M80").Value = -1445.90905
$ws.Range("N80").Value = -5146

$ws.Range("H83").Value = 865.0714
$ws.Range("I83").Value = 814.63635
$ws.Range("J83").Value = 1050
$ws.Range("K83").Value = 7331.72715
$ws.Range("L83").Value = 9450
$ws.Range("M83").Value = -2339.72715
$ws.Range("N83").Value = -19434

$ws.Range("H88").Value = 1431.2858
$ws.Range("I88").Value = 728.7273
$ws.Range("J88").Value = 1753.2916
$ws.Range("K88").Value = 728.7273
$ws.Range("L88").Value = 1753.2916
$ws.Range("M88").Value = -322.7273
$ws.Range("N88").Value = -2565.2916

$ws.Range("H91").Value = 1431.2858
$ws.Range("I91").Value = 728.7273
$ws.Range("J91").Value = 1753.2916
$ws.Range("K91").Value = 728.7273
$ws.Range("L91").Value = 1753.2916
$ws.Range("M91").Value = 675.2727
$ws.Range("N91").Value = -4561.2916

$ws.Range("H98").Value = 4911.6665
$ws.Range("I98").Value = 4911.6665
$ws.Range("K98").Value = 4911.6665
$ws.Range("M98").Value = -3413.6665

$ws.Range("H122").Value = 4911.6665
$ws.Range("I122").Value = 4911.6665
$ws.Range("K122").Value = 14734.9995
$ws.Range("M122").Value = -12284.9995

$ws.Range("H137").Value = 7956.647
$ws.Range("I137").Value = 11716.4
$ws.Range("J137").Value = 2585.5715
$ws.Range("K137").Value = 35149.2
$ws.Range("L137").Value = 7756.7145
$ws.Range("M137").Value = -32599.2
$ws.Range("N137").Value = -12856.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 2489.6
$ws.Range("I25").Value = 1112
$ws.Range("J25").Value = 8000
$ws.Range("K25").Value = 1112
$ws.Range("L25").Value = 8000
$ws.Range("M25").Value = -710
$ws.Range("N25").Value = -8804

$ws.Range("H110").Value = 1077.091
$ws.Range("I110").Value = 882.1177
$ws.Range("J110").Value = 1740
$ws.Range("K110").Value = 882.1177
$ws.Range("L110").Value = 1740
$ws.Range("M110").Value = 1162.8823
$ws.Range("N110").Value = -5830

$ws.Range("H122").Value = 27780050
$ws.Range("I122").Value = 62501800
$ws.Range("J122").Value = 2649.9
$ws.Range("K122").Value = 187505400
$ws.Range("L122").Value = 7949.700000000001
$ws.Range("M122").Value = -187502950
$ws.Range("N122").Value = -12849.7

$ws.Range("H132").Value = 2784.318
$ws.Range("I132").Value = 2780.8333
$ws.Range("J132").Value = 2788.5
$ws.Range("K132").Value = 8342.499899999999
$ws.Range("L132").Value = 8365.5
$ws.Range("M132").Value = -5812.499899999999
$ws.Range("N132").Value = -13425.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2246.8096
$ws.Range("I31").Value = 1834
$ws.Range("J31").Value = 2500.8462
$ws.Range("K31").Value = 1834
$ws.Range("L31").Value = 2500.8462
$ws.Range("M31").Value = -1539
$ws.Range("N31").Value = -3090.8462

$ws.Range("H34").Value = 2246.8096
$ws.Range("I34").Value = 1834
$ws.Range("J34").Value = 2500.8462
$ws.Range("K34").Value = 1834
$ws.Range("L34").Value = 2500.8462
$ws.Range("M34").Value = -1632
$ws.Range("N34").Value = -2904.8462

$ws.Range("H60").Value = 8964.333000000001
$ws.Range("I60").Value = 2093
$ws.Range("J60").Value = 12400
$ws.Range("K60").Value = 2093
$ws.Range("L60").Value = 12400
$ws.Range("M60").Value = -1582
$ws.Range("N60").Value = -13422

$ws.Range("H74").Value = 19800
$ws.Range("J74").Value = 19800
$ws.Range("L74").Value = 19800
$ws.Range("N74").Value = -21548

$ws.Range("H77").Value = 19800
$ws.Range("J77").Value = 19800
$ws.Range("L77").Value = 59400
$ws.Range("N77").Value = -68136

$ws.Range("H86").Value = 2640
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 2830
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 2830
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -5076

$ws.Range("H89").Value = 2640
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 2830
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 14150
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -25382

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 142857630
$ws.Range("I46").Value = 350
$ws.Range("J46").Value = 333334000
$ws.Range("K46").Value = 1050
$ws.Range("L46").Value = 1000002000
$ws.Range("M46").Value = -959
$ws.Range("N46").Value = -1000002182

$ws.Range("H68").Value = 940.1707
$ws.Range("I68").Value = 611.73334
$ws.Range("J68").Value = 1129.6538
$ws.Range("K68").Value = 1835.20002
$ws.Range("L68").Value = 3388.9614
$ws.Range("M68").Value = -1024.20002
$ws.Range("N68").Value = -5010.9614

$ws.Range("H71").Value = 940.1707
$ws.Range("I71").Value = 611.73334
$ws.Range("J71").Value = 1129.6538
$ws.Range("K71").Value = 5505.60006
$ws.Range("L71").Value = 10166.8842
$ws.Range("M71").Value = -1449.60006
$ws.Range("N71").Value = -18278.8842

$ws.Range("H131").Value = 1236267.6
$ws.Range("J131").Value = 1317450.6
$ws.Range("L131").Value = 3952351.8
$ws.Range("N131").Value = -3962431.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 174000.86
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 202501
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 202501
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -204497

$ws.Range("H83").Value = 174000.86
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 202501
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 1012505
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -1022489

$ws.Range("H113").Value = 2516.2727
$ws.Range("I113").Value = 3275.8
$ws.Range("J113").Value = 1883.3334
$ws.Range("K113").Value = 3275.8
$ws.Range("L113").Value = 1883.3334
$ws.Range("M113").Value = -1105.8
$ws.Range("N113").Value = -6223.3334

$ws.Range("H126").Value = 1657.5
$ws.Range("I126").Value = 1365
$ws.Range("J126").Value = 1950
$ws.Range("K126").Value = 4095
$ws.Range("L126").Value = 5850
$ws.Range("M126").Value = -1625
$ws.Range("N126").Value = -10790

$ws.Range("H132").Value = 1610.0358
$ws.Range("I132").Value = 1297.6316
$ws.Range("J132").Value = 2269.5557
$ws.Range("K132").Value = 3892.8948
$ws.Range("L132").Value = 6808.6671
$ws.Range("M132").Value = -1362.8948
$ws.Range("N132").Value = -11868.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1858.2858
$ws.Range("I7").Value = 1858.2858
$ws.Range("K7").Value = 1858.2858
$ws.Range("M7").Value = -1746.2858

$ws.Range("H68").Value = 2912.3333
$ws.Range("I68").Value = 2423.5454
$ws.Range("J68").Value = 3450
$ws.Range("K68").Value = 2423.5454
$ws.Range("L68").Value = 3450
$ws.Range("M68").Value = -1674.5454
$ws.Range("N68").Value = -4948

$ws.Range("H71").Value = 2912.3333
$ws.Range("I71").Value = 2423.5454
$ws.Range("J71").Value = 3450
$ws.Range("K71").Value = 12117.727
$ws.Range("L71").Value = 17250
$ws.Range("M71").Value = -8373.726999999999
$ws.Range("N71").Value = -24738

$ws.Range("H82").Value = 2498.5715
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 2581.6667
$ws.Range("K82").Value = 2000
$ws.Range("L82").Value = 2581.6667
$ws.Range("M82").Value = -1639
$ws.Range("N82").Value = -3303.6667

$ws.Range("H85").Value = 2498.5715
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 2581.6667
$ws.Range("K85").Value = 2000
$ws.Range("L85").Value = 2581.6667
$ws.Range("M85").Value = -752
$ws.Range("N85").Value = -5077.6667

$ws.Range("H126").Value = 1858.2858
$ws.Range("I126").Value = 1858.2858
$ws.Range("K126").Value = 5574.857400000001
$ws.Range("M126").Value = -3104.857400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 479.57144
$ws.Range("I113").Value = 459
$ws.Range("J113").Value = 603
$ws.Range("K113").Value = 1377
$ws.Range("L113").Value = 1809
$ws.Range("M113").Value = 793
$ws.Range("N113").Value = -6149
